$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 70 (shifts old rows 70-126 down to 72-128)
$ws.Range("A70:A71").EntireRow.Insert()

# New row 70: same dimension/category as the (now shifted) row 72, with updated date + volume
$ws.Cells.Item(70, 1).Value = 8
$ws.Cells.Item(70, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44981
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112027
$ws.Cells.Item(70, 7).Value = "Melón"
$ws.Cells.Item(70, 8).Value = "Tuna"
$ws.Cells.Item(70, 9).Value = "Extra"
$ws.Cells.Item(70, 10).Value = 2000
$ws.Cells.Item(70, 11).Value = 1400
$ws.Cells.Item(70, 12).Value = 1500
$ws.Cells.Item(70, 13).Value = 1450
$ws.Cells.Item(70, 14).Value = "`$/unidad"
$ws.Cells.Item(70, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 16).Value = 1450
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(70, 18).Value = "Hortaliza"

# New row 71
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44981
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112027
$ws.Cells.Item(71, 7).Value = "Melón"
$ws.Cells.Item(71, 8).Value = "Tuna"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 1600
$ws.Cells.Item(71, 11).Value = 1100
$ws.Cells.Item(71, 12).Value = 1200
$ws.Cells.Item(71, 13).Value = 1150
$ws.Cells.Item(71, 14).Value = "`$/unidad"
$ws.Cells.Item(71, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 16).Value = 1150
$ws.Cells.Item(71, 17).Value = 1
$ws.Cells.Item(71, 18).Value = "Hortaliza"
